$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove existing hyperlink tracking (will re-add after the shift)
$ws.Hyperlinks.Delete()

# 2. Delete old row 3 (orig "2/2022/Seminarios" row with no link) - shifts rows 4-21 up to 3-20
$ws.Rows("3:3").Delete()

# 3. Add first new "Seminarios" row (row 21) - re-introduces the "Seminarios" shared string
#    and introduces the first brand-new shared string (the YouTube link) right away,
#    matching the original authoring order captured in the diff.
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = 2022
$ws.Cells.Item(21, 3).Value = 'Seminarios'
$ws.Cells.Item(21, 4).Value = 'https://www.youtube.com/watch?v=VSUEi50tkAI'

# 4. Rename "Nota ecoinformatica" -> "Notas ecoinformaticas" for rows 4-20
for ($r = 4; $r -le 20; $r++) {
  $ws.Cells.Item($r, 3).Value = 'Notas ecoinformáticas'
}

# 5. Add remaining new "Seminarios" rows (22-36)
$ws.Cells.Item(22, 1).Value = 2
$ws.Cells.Item(22, 2).Value = 2022
$ws.Cells.Item(22, 3).Value = 'Seminarios'
$ws.Cells.Item(22, 4).Value = 'https://www.youtube.com/watch?v=usB7reMJxLU'
$ws.Cells.Item(23, 1).Value = 3
$ws.Cells.Item(23, 2).Value = 2022
$ws.Cells.Item(23, 3).Value = 'Seminarios'
$ws.Cells.Item(23, 4).Value = 'https://www.youtube.com/watch?v=OKEtldANpHA'
$ws.Cells.Item(24, 1).Value = 5
$ws.Cells.Item(24, 2).Value = 2022
$ws.Cells.Item(24, 3).Value = 'Seminarios'
$ws.Cells.Item(24, 4).Value = 'https://www.youtube.com/watch?v=niPkyIrEv-k'
$ws.Cells.Item(25, 1).Value = 6
$ws.Cells.Item(25, 2).Value = 2022
$ws.Cells.Item(25, 3).Value = 'Seminarios'
$ws.Cells.Item(25, 4).Value = 'https://www.youtube.com/watch?v=oAC7DVWAMRc'
$ws.Cells.Item(26, 1).Value = 7
$ws.Cells.Item(26, 2).Value = 2022
$ws.Cells.Item(26, 3).Value = 'Seminarios'
$ws.Cells.Item(26, 4).Value = 'https://www.youtube.com/watch?v=0_73g05Wjgc'
$ws.Cells.Item(27, 1).Value = 9
$ws.Cells.Item(27, 2).Value = 2022
$ws.Cells.Item(27, 3).Value = 'Seminarios'
$ws.Cells.Item(27, 4).Value = 'https://www.youtube.com/watch?v=ybD4GM_OB3M'
$ws.Cells.Item(28, 1).Value = 10
$ws.Cells.Item(28, 2).Value = 2022
$ws.Cells.Item(28, 3).Value = 'Seminarios'
$ws.Cells.Item(28, 4).Value = 'https://www.youtube.com/watch?v=1ye4v3ugAfc'
$ws.Cells.Item(29, 1).Value = 11
$ws.Cells.Item(29, 2).Value = 2022
$ws.Cells.Item(29, 3).Value = 'Seminarios'
$ws.Cells.Item(29, 4).Value = 'https://www.youtube.com/watch?v=ZwKuy6E-GhE'
$ws.Cells.Item(30, 1).Value = 2
$ws.Cells.Item(30, 2).Value = 2023
$ws.Cells.Item(30, 3).Value = 'Seminarios'
$ws.Cells.Item(30, 4).Value = 'https://www.youtube.com/watch?v=nkEM6ny_E9U'
$ws.Cells.Item(31, 1).Value = 5
$ws.Cells.Item(31, 2).Value = 2023
$ws.Cells.Item(31, 3).Value = 'Seminarios'
$ws.Cells.Item(31, 4).Value = 'https://www.youtube.com/watch?v=C9Gyah5XES0'
$ws.Cells.Item(32, 1).Value = 9
$ws.Cells.Item(32, 2).Value = 2023
$ws.Cells.Item(32, 3).Value = 'Seminarios'
$ws.Cells.Item(32, 4).Value = 'https://www.youtube.com/watch?v=Fy5YNSe-btA'
$ws.Cells.Item(33, 1).Value = 1
$ws.Cells.Item(33, 2).Value = 2024
$ws.Cells.Item(33, 3).Value = 'Seminarios'
$ws.Cells.Item(33, 4).Value = 'https://www.youtube.com/watch?v=CPimYLdDI5o'
$ws.Cells.Item(34, 1).Value = 3
$ws.Cells.Item(34, 2).Value = 2024
$ws.Cells.Item(34, 3).Value = 'Seminarios'
$ws.Cells.Item(34, 4).Value = 'https://www.youtube.com/watch?v=vzjqLjKOPPM'
$ws.Cells.Item(35, 1).Value = 4
$ws.Cells.Item(35, 2).Value = 2024
$ws.Cells.Item(35, 3).Value = 'Seminarios'
$ws.Cells.Item(35, 4).Value = 'https://www.youtube.com/watch?v=KLN2wMnivVA'
$ws.Cells.Item(36, 1).Value = 5
$ws.Cells.Item(36, 2).Value = 2024
$ws.Cells.Item(36, 3).Value = 'Seminarios'
$ws.Cells.Item(36, 4).Value = 'https://www.youtube.com/watch?v=txJfuPeodQw'

# 6. Re-add the hyperlink, now anchored at D4 (was D5 before the row shift).
#    D4 already inherited the correct "Hipervinculo" cell style (s="1") from the
#    original D5 when the row above it was deleted, so stash that style first and
#    reapply it afterward - Hyperlinks.Add() otherwise stamps the cell with a
#    freshly duplicated (but visually identical) style entry.
$goodStyle = $ws.Range("D4").Style
$ws.Hyperlinks.Add($ws.Range("D4"), "https://doi.org/10.7818/ECOS.2017.26-1.20")
$ws.Range("D4").Style = $goodStyle

# 7. Match final cursor/selection position from the authored workbook
$ws.Range("A37").Select()
